# Updated small form factor headset 3D prints
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Small form factor headset")

# Make this sheet the active/selected tab (previously it was "Monocular eyepiece")
$ws.Activate()

# Row 25: the old combined enclosure print (V1) becomes the V2 enclosure "back" print
$ws.Range("A25").Value = "headset enclosure back"
$ws.Range("B25").Value = "miniHeadset_V2_backing.stl"

# Insert a new row right after row 25 for the new "headset eyepiece" 3D print
# (this shifts the old "lens holder" / "display holder" rows down by one)
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "headset eyepiece"
$ws.Range("B26").Value = "miniHeadset_V2_eyepiece.stl"
$ws.Range("C26").Value = "cm"
$ws.Range("D26").Value = "see README for 3D printing instructions"
$ws.Range("E26").Value = "~$5-10 + shipping"
$ws.Range("F26").Value = 2

# Leave the selection where the author left it when saving
$null = $ws.Range("B32").Select()
